# Update betting-odds values in Sheet1 to match the 2025-05-20 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Crystal Palace vs Wolves)
$ws.Range("U2").Value = 13
$ws.Range("AA2").Value = 7
$ws.Range("AI2").Value = 26

# Row 3 (Manchester City vs Bournemouth)
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 19
$ws.Range("Z3").Value = 19
$ws.Range("AB3").Value = 17

# Row 4 (Den Bosch vs Telstar)
$ws.Range("G4").Value = 2.9
$ws.Range("I4").Value = 2.35
$ws.Range("J4").Value = 1.06
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 1.33
$ws.Range("M4").Value = 3.25
$ws.Range("N4").Value = 2.05
$ws.Range("O4").Value = 1.75
$ws.Range("P4").Value = 1.4
$ws.Range("Q4").Value = 2.75
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.95
$ws.Range("T4").Value = 8.5
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 34
$ws.Range("Z4").Value = 9.5
$ws.Range("AB4").Value = 15
$ws.Range("AC4").Value = 51
$ws.Range("AD4").Value = 301
$ws.Range("AE4").Value = 7.5
$ws.Range("AF4").Value = 11
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 29

# Row 8 (Hebar vs Botev Vratsa)
$ws.Range("G8").Value = 7.5
$ws.Range("H8").Value = 4.1
$ws.Range("I8").Value = 1.48
$ws.Range("T8").Value = 17
$ws.Range("W8").Value = 81
$ws.Range("AB8").Value = 19

# Row 10 (CSKA 1948 Sofia II vs CSKA Sofia II) - odds were missing, now populated
$ws.Range("G10").Value = 2.47
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 2.47
$ws.Range("N10").Value = 1.91
$ws.Range("O10").Value = 1.7
$ws.Range("P10").Value = 1.39
$ws.Range("Q10").Value = 2.42
$ws.Range("T10").Value = 6.9
$ws.Range("U10").Value = 10
$ws.Range("V10").Value = 8.25
$ws.Range("W10").Value = 21
$ws.Range("X10").Value = 17
$ws.Range("Y10").Value = 25
$ws.Range("Z10").Value = 9.25
$ws.Range("AA10").Value = 5.6
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 50
$ws.Range("AD10").Value = 350
$ws.Range("AE10").Value = 6.9
$ws.Range("AF10").Value = 10
$ws.Range("AG10").Value = 8.25
$ws.Range("AH10").Value = 21
$ws.Range("AI10").Value = 17
$ws.Range("AJ10").Value = 25

# Row 24
$ws.Range("N24").Value = 1.53
$ws.Range("O24").Value = 2.38

# Row 26
$ws.Range("G26").Value = 2.1
$ws.Range("I26").Value = 2.88
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 1.03
$ws.Range("W26").Value = 21
$ws.Range("AC26").Value = 29
$ws.Range("AH26").Value = 34
$ws.Range("AI26").Value = 21

# Row 28
$ws.Range("G28").Value = 1.5
$ws.Range("H28").Value = 4.2
$ws.Range("I28").Value = 6
$ws.Range("L28").Value = 1.22
$ws.Range("M28").Value = 4
$ws.Range("N28").Value = 1.65
$ws.Range("O28").Value = 2.15
$ws.Range("P28").Value = 1.34
$ws.Range("Q28").Value = 3.1
$ws.Range("R28").Value = 1.78
$ws.Range("S28").Value = 1.93
$ws.Range("T28").Value = 7.2
$ws.Range("U28").Value = 7.9
$ws.Range("W28").Value = 11.25
$ws.Range("X28").Value = 12.5
$ws.Range("Y28").Value = 25
$ws.Range("AA28").Value = 8.75
$ws.Range("AB28").Value = 17.5
$ws.Range("AC28").Value = 75
$ws.Range("AD28").Value = 500
$ws.Range("AE28").Value = 16.5
$ws.Range("AF28").Value = 45
$ws.Range("AG28").Value = 19.5
$ws.Range("AH28").Value = 150
$ws.Range("AI28").Value = 60
$ws.Range("AJ28").Value = 55
